$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 3529.1667
$ws.Range("I2").Value = 3529.1667
$ws.Range("K2").Value = 3529.1667
$ws.Range("M2").Value = -3416.1667

# Row 70
$ws.Range("H70").Value = 2509.6
$ws.Range("J70").Value = 3000.6
$ws.Range("L70").Value = 9001.799999999999
$ws.Range("N70").Value = -9541.799999999999

# Row 73
$ws.Range("H73").Value = 2509.6
$ws.Range("J73").Value = 3000.6
$ws.Range("L73").Value = 9001.799999999999
$ws.Range("N73").Value = -10873.8

# Row 115
$ws.Range("H115").Value = 2475
$ws.Range("J115").Value = 3900
$ws.Range("L115").Value = 11700
$ws.Range("N115").Value = -14834

# Row 137
$ws.Range("H137").Value = 38333.82
$ws.Range("I137").Value = 51506.168
$ws.Range("K137").Value = 154518.504
$ws.Range("M137").Value = -151968.504

# Row 138
$ws.Range("H138").Value = 3650.8147
$ws.Range("I138").Value = 2675.8
$ws.Range("J138").Value = 4025.8206
$ws.Range("K138").Value = 8027.400000000001
$ws.Range("L138").Value = 12077.4618
$ws.Range("M138").Value = -2887.400000000001
$ws.Range("N138").Value = -22357.4618

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 17767.092
$ws.Range("I32").Value = 15264.625
$ws.Range("J32").Value = 24440.334
$ws.Range("K32").Value = 15264.625
$ws.Range("L32").Value = 24440.334
$ws.Range("M32").Value = -14977.625
$ws.Range("N32").Value = -25014.334

# Row 74
$ws.Range("H74").Value = 42355.39
$ws.Range("I74").Value = 3047.3845
$ws.Range("K74").Value = 3047.3845
$ws.Range("M74").Value = -2173.3845

# Row 77
$ws.Range("H77").Value = 42355.39
$ws.Range("I77").Value = 3047.3845
$ws.Range("K77").Value = 15236.9225
$ws.Range("M77").Value = -10868.9225

$ws = $wb.Worksheets.Item("BSM")
# Row 81
$ws.Range("H81").Value = 10246.75
$ws.Range("J81").Value = 10246.75
$ws.Range("L81").Value = 10246.75
$ws.Range("N81").Value = -12368.75

# Row 84
$ws.Range("H84").Value = 10246.75
$ws.Range("J84").Value = 10246.75
$ws.Range("L84").Value = 30740.25
$ws.Range("N84").Value = -41348.25

# Row 94
$ws.Range("H94").Value = 3476006
$ws.Range("I94").Value = 5001613.5
$ws.Range("J94").Value = 8716.817999999999
$ws.Range("K94").Value = 5001613.5
$ws.Range("L94").Value = 8716.817999999999
$ws.Range("M94").Value = -5001162.5
$ws.Range("N94").Value = -9618.817999999999

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 884.4211
$ws.Range("I16").Value = 530.25
$ws.Range("K16").Value = 530.25
$ws.Range("M16").Value = -243.25

# Row 58
$ws.Range("H58").Value = 10856.923
$ws.Range("I58").Value = 13294.444
$ws.Range("J58").Value = 5372.5
$ws.Range("K58").Value = 13294.444
$ws.Range("L58").Value = 5372.5
$ws.Range("M58").Value = -13091.444
$ws.Range("N58").Value = -5778.5

# Row 60
$ws.Range("H60").Value = 27449
$ws.Range("J60").Value = 70000
$ws.Range("L60").Value = 70000
$ws.Range("N60").Value = -71022

# Row 99
$ws.Range("H99").Value = 4183.1333
$ws.Range("I99").Value = 4008.889
$ws.Range("J99").Value = 4444.5
$ws.Range("K99").Value = 4008.889
$ws.Range("L99").Value = 4444.5
$ws.Range("M99").Value = -2510.889
$ws.Range("N99").Value = -7440.5

# Row 113
$ws.Range("H113").Value = 884.4211
$ws.Range("I113").Value = 530.25
$ws.Range("K113").Value = 530.25
$ws.Range("M113").Value = 1639.75

# Row 126
$ws.Range("H126").Value = 4183.1333
$ws.Range("I126").Value = 4008.889
$ws.Range("J126").Value = 4444.5
$ws.Range("K126").Value = 12026.667
$ws.Range("L126").Value = 13333.5
$ws.Range("M126").Value = -9556.667000000001
$ws.Range("N126").Value = -18273.5

# Row 132
$ws.Range("H132").Value = 44105.066
$ws.Range("I132").Value = 31370.03
$ws.Range("K132").Value = 94110.09
$ws.Range("M132").Value = -91580.09

# Row 136
$ws.Range("H136").Value = 10856.923
$ws.Range("I136").Value = 13294.444
$ws.Range("J136").Value = 5372.5
$ws.Range("K136").Value = 39883.33199999999
$ws.Range("L136").Value = 16117.5
$ws.Range("M136").Value = -37333.33199999999
$ws.Range("N136").Value = -21217.5

# Row 141
$ws.Range("H141").Value = 166870.1
$ws.Range("J141").Value = 166870.1
$ws.Range("L141").Value = 166870.1
$ws.Range("N141").Value = -177230.1

$ws = $wb.Worksheets.Item("CUL")
# Row 26
$ws.Range("H26").Value = 161.25
$ws.Range("I26").Value = 161.25
$ws.Range("K26").Value = 483.75
$ws.Range("M26").Value = -195.75

# Row 129
$ws.Range("H129").Value = 1539836
$ws.Range("J129").Value = 1765
$ws.Range("L129").Value = 5295
$ws.Range("N129").Value = -15295

# Row 137
$ws.Range("H137").Value = 7423.294
$ws.Range("J137").Value = 7785.5713
$ws.Range("L137").Value = 23356.7139
$ws.Range("N137").Value = -33556.7139

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 1834156.4
$ws.Range("I80").Value = 2278875.2
$ws.Range("J80").Value = 499999.8
$ws.Range("K80").Value = 2278875.2
$ws.Range("L80").Value = 499999.8
$ws.Range("M80").Value = -2277877.2
$ws.Range("N80").Value = -501995.8

# Row 83
$ws.Range("H83").Value = 1834156.4
$ws.Range("I83").Value = 2278875.2
$ws.Range("J83").Value = 499999.8
$ws.Range("K83").Value = 11394376
$ws.Range("L83").Value = 2499999
$ws.Range("M83").Value = -11389384
$ws.Range("N83").Value = -2509983

# Row 102
$ws.Range("H102").Value = 6427292.5
$ws.Range("I102").Value = 12347184
$ws.Range("K102").Value = 12347184
$ws.Range("M102").Value = -12345562

# Row 136
$ws.Range("H136").Value = 57996.332
$ws.Range("J136").Value = 57996.332
$ws.Range("L136").Value = 173988.996
$ws.Range("N136").Value = -179088.996

$ws = $wb.Worksheets.Item("LTW")
# Row 2
$ws.Range("H2").Value = 2008499.8
$ws.Range("I2").Value = 5001500
$ws.Range("K2").Value = 5001500
$ws.Range("M2").Value = -5001388

# Row 40
$ws.Range("H40").Value = 7410.579
$ws.Range("I40").Value = 4199.4443
$ws.Range("J40").Value = 10300.6
$ws.Range("K40").Value = 4199.4443
$ws.Range("L40").Value = 10300.6
$ws.Range("M40").Value = -4063.4443
$ws.Range("N40").Value = -10572.6

# Row 46
$ws.Range("H46").Value = 4265.125
$ws.Range("I46").Value = 922.1111
$ws.Range("K46").Value = 922.1111
$ws.Range("M46").Value = -734.1111

# Row 93
$ws.Range("H93").Value = 33348042
$ws.Range("I93").Value = 55559736
$ws.Range("K93").Value = 55559736
$ws.Range("M93").Value = -55558488

# Row 94
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").Value = $null

# Row 132
$ws.Range("H132").Value = 15435.077
$ws.Range("I132").Value = 15435.077
$ws.Range("K132").Value = 46305.231
$ws.Range("M132").Value = -43775.231

$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 6672178.5
$ws.Range("I81").Value = 9263204
$ws.Range("J81").Value = 9541.286
$ws.Range("K81").Value = 18526408
$ws.Range("L81").Value = 19082.572
$ws.Range("M81").Value = -18525347
$ws.Range("N81").Value = -21204.572

# Row 84
$ws.Range("H84").Value = 6672178.5
$ws.Range("I84").Value = 9263204
$ws.Range("J84").Value = 9541.286
$ws.Range("K84").Value = 92632040
$ws.Range("L84").Value = 95412.86
$ws.Range("M84").Value = -92626736
$ws.Range("N84").Value = -106020.86

# Row 115
$ws.Range("H115").Value = 46000
$ws.Range("J115").Value = 46000
$ws.Range("L115").Value = 46000
$ws.Range("N115").Value = -49134

# Row 132
$ws.Range("H132").Value = 18383320
$ws.Range("I132").Value = 19234400
$ws.Range("J132").Value = 3631295.8
$ws.Range("K132").Value = 57703200
$ws.Range("L132").Value = 10893887.4
$ws.Range("M132").Value = -57700670
$ws.Range("N132").Value = -10898947.4

# Row 136
$ws.Range("H136").Value = 3034.9185
$ws.Range("I136").Value = 2748.7693
$ws.Range("K136").Value = 8246.3079
$ws.Range("M136").Value = -5696.3079
